$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Range("C4").Value2 = 264
Write-Host ("New Value2: " + $ws.Range("C4").Value2)
Write-Host ("L3 formula: " + $ws.Range("L3").Formula)
Write-Host ("L3 Value2: " + $ws.Range("L3").Value2)
Write-Host ("L4 Value2: " + $ws.Range("L4").Value2)
Write-Host ("L5 Value2: " + $ws.Range("L5").Value2)
